# Refresh cryptos list values (prices / 1h volume %) per scheduled GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.543.18'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '2.246.27'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = '''305.85'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').Value = '''94.63'
$ws.Range('E6').Value = '  -1.55%  '
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('D10').Value = '''34.74'
$ws.Range('E10').Value = '  -0.62%  '
$ws.Range('E11').Value = '  -1.51%  '
$ws.Range('E13').Value = '  -0.04%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '2.386.20'
$ws.Range('E14').Value = '  +1.91%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.588.65'
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('E17').Value = '  -0.35%  '
$ws.Range('D18').Value = '44.337.19'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('E19').Value = '  -2.98%  '
$ws.Range('D20').Value = '''6.19'
$ws.Range('E20').Value = '  -2.81%  '
$ws.Range('D21').Value = '''11.73'
$ws.Range('E21').Value = '  -3.24%  '
$ws.Range('D22').Value = '''65.23'
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('D23').Value = '''236.69'
$ws.Range('E23').Value = '  -0.66%  '
$ws.Range('D24').Value = '''2.95'
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('E25').Value = '  -1.58%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').Value = '''2.33'
$ws.Range('E27').Value = '  +5.51%  '
$ws.Range('D28').Value = '''9.77'
$ws.Range('E29').Value = '  -3.97%  '
$ws.Range('D30').Value = '''5.90'
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').Value = '''19.95'
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('D32').Value = '''149.40'
$ws.Range('E32').Value = '  -2.21%  '
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('D35').Value = '''3.21'
$ws.Range('E35').Value = '  -1.18%  '
$ws.Range('D36').Value = '''0.109'
$ws.Range('E37').Value = '  -1.21%  '
$ws.Range('D38').Value = '''1.87'
$ws.Range('E38').Value = '  +5.94%  '
$ws.Range('D39').Value = '''15.22'
$ws.Range('E39').Value = '  +4.90%  '
$ws.Range('E40').Value = '  -6.01%  '
$ws.Range('E41').Value = '  -1.97%  '
$ws.Range('E42').Value = '  -0.89%  '
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('D44').Value = '1.802.65'
$ws.Range('E44').Value = '  +2.76%  '
$ws.Range('E45').Value = '  +11.28%  '
$ws.Range('D46').Value = '''82.07'
$ws.Range('E46').Value = '  -1.19%  '
$ws.Range('E47').Value = '  -2.44%  '
$ws.Range('D48').Value = '''98.43'
$ws.Range('E48').Value = '  -1.73%  '
$ws.Range('E49').Value = '  -2.50%  '
$ws.Range('D50').Value = '''68.74'
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('D51').Value = '''53.80'
$ws.Range('E51').Value = '  -1.84%  '
